$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NAICS-style investment codes (one per data row, E2:E192)
$naicsCodes = @(
    "99100000",
    "99110000",
    "99111000",
    "99111100",
    "99111200",
    "99111210",
    "99111211",
    "99111212",
    "99111220",
    "99111300",
    "99111400",
    "99111500",
    "99111600",
    "99112000",
    "99113000",
    "99113100",
    "99113110",
    "99113111",
    "99113112",
    "99113120",
    "99113200",
    "99114000",
    "99114100",
    "99114200",
    "99115000",
    "99115100",
    "99115200",
    "99115300",
    "99115400",
    "99115600",
    "99115610",
    "99115620",
    "99115700",
    "99115800",
    "99115900",
    "99115999",
    "99116000",
    "99116100",
    "99116110",
    "99116120",
    "99116200",
    "99116210",
    "99116220",
    "99116230",
    "99116240",
    "99116250",
    "99116300",
    "99116310",
    "99116320",
    "99200000",
    "99210000",
    "99211000",
    "99211100",
    "99211200",
    "99211210",
    "99211211",
    "99211212",
    "99211220",
    "99211300",
    "99211400",
    "99211500",
    "99211600",
    "99212000",
    "99213000",
    "99213100",
    "99213110",
    "99213111",
    "99213112",
    "99213120",
    "99213200",
    "99214000",
    "99214100",
    "99214200",
    "99215000",
    "99215100",
    "99215200",
    "99215300",
    "99215400",
    "99215500",
    "99215510",
    "99215520",
    "99215600",
    "99215700",
    "99215800",
    "99215900",
    "99220000",
    "99221000",
    "99221100",
    "99221200",
    "99222000",
    "99222100",
    "99222200",
    "99222300",
    "99222400",
    "99222500",
    "99230000",
    "99231000",
    "99231100",
    "99231200",
    "99300000",
    "99310000",
    "99311000",
    "99311100",
    "99311200",
    "99311300",
    "99311310",
    "99311320",
    "99311400",
    "99311500",
    "99311600",
    "99312000",
    "99312100",
    "99312200",
    "99312210",
    "99312220",
    "99312300",
    "99312400",
    "99312500",
    "99312600",
    "99313000",
    "99313100",
    "99313110",
    "99313120",
    "99313200",
    "99313300",
    "99313400",
    "99313500",
    "99314000",
    "99314100",
    "99314110",
    "99314120",
    "99314200",
    "99314210",
    "99314220",
    "99314300",
    "99314310",
    "99314320",
    "99314400",
    "99314500",
    "99314600",
    "99314610",
    "99314620",
    "99314700",
    "99315000",
    "99320000",
    "99400000",
    "99410000",
    "99411000",
    "99411100",
    "99411200",
    "99411300",
    "99411310",
    "99411320",
    "99411400",
    "99411500",
    "99411600",
    "99412000",
    "99412100",
    "99412200",
    "99412210",
    "99412220",
    "99412300",
    "99412400",
    "99412500",
    "99412600",
    "99413000",
    "99413100",
    "99413110",
    "99413120",
    "99413200",
    "99413300",
    "99413400",
    "99413500",
    "99414000",
    "99414100",
    "99414110",
    "99414120",
    "99414200",
    "99414210",
    "99414220",
    "99414300",
    "99414310",
    "99414320",
    "99414400",
    "99414500",
    "99414600",
    "99414610",
    "99414620",
    "99414700",
    "99415000",
    "99420000"
)

# Write each code as a text formula first so it is not coerced into a
# number by general cell formatting, then collapse the formulas down to
# plain literal text values via copy / paste-special (values only).
for ($i = 0; $i -lt $naicsCodes.Length; $i++) {
    $row = $i + 2
    $code = $naicsCodes[$i]
    $ws.Range("E$row").Formula = '="' + $code + '"'
}

$dataRange = $ws.Range("E2:E192")
$dataRange.Copy() | Out-Null
$dataRange.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Restore the view to the top of the sheet with E2 selected (matches the
# author's saved selection / scroll position after the edit).
$ws.Range("E2").Select() | Out-Null
